$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new weekly date columns right after column A (before the old column B),
# shifting all existing week columns (old B..V) to the right by 9 (new K..AE).
$ws.Range("B1:J1").EntireColumn.Insert()

# Match the existing column width formatting (8.0 chars) across the now-wider used range.
$ws.Range("C1:AE1").ColumnWidth = 7.1

# New header row: 9 additional weekly date labels, newest-first, continuing the series
# that used to start at Jun_09 (now pushed out to column K).
$ws.Cells.Item(1, 2).Value = "Sep_08"
$ws.Cells.Item(1, 3).Value = "Aug_25"
$ws.Cells.Item(1, 4).Value = "Aug_04"
$ws.Cells.Item(1, 5).Value = "Jul_23"
$ws.Cells.Item(1, 6).Value = "Jul_17"
$ws.Cells.Item(1, 7).Value = "Jul_07"
$ws.Cells.Item(1, 8).Value = "Jun_30"
$ws.Cells.Item(1, 9).Value = "Jun_24"
$ws.Cells.Item(1, 10).Value = "Jun_16"

# Fill the new weekly columns (B..J) for every analyst/firm data row with the default
# "UN" (unchanged) rating marker, same as all the other pre-existing weekly columns.
for ($r = 2; $r -le 33; $r++) {
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

# New rating-change note for Evercore ISI (row 29) on the week of 8/12/2019 -> Aug_25 column.
$ws.Cells.Item(29, 3).Value = "8/12/2019,Reiterates,Sell,"
